$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 42.517849
$ws.Cells.Item(2, 8).Value = 127.553547
$ws.Cells.Item(2, 9).Value = 0.02311569285614191
$ws.Cells.Item(2, 10).Value = 0.02311569285614191
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 4.220261333333333
$ws.Cells.Item(2, 14).Value = 12.660784
$ws.Cells.Item(2, 15).Value = 0.6739259863235564
$ws.Cells.Item(2, 16).Value = 0.6739259863235564
$ws.Cells.Item(2, 17).Value = 179.4364341112053
$ws.Cells.Item(2, 18).Value = 1614.927907000848
$ws.Cells.Item(2, 19).Value = 0.01557826610762783
$ws.Cells.Item(2, 20).Value = 0.01557826610762782

$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 42.517849
$ws.Cells.Item(3, 8).Value = 127.553547
$ws.Cells.Item(3, 9).Value = 0.02311569285614191
$ws.Cells.Item(3, 10).Value = 0.02311569285614191
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.296447666666667
$ws.Cells.Item(3, 14).Value = 3.889343
$ws.Cells.Item(3, 15).Value = 0.2070274097896007
$ws.Cells.Item(3, 16).Value = 0.2070274097896007
$ws.Cells.Item(3, 17).Value = 55.12216612773567
$ws.Cells.Item(3, 18).Value = 496.099495149621
$ws.Cells.Item(3, 19).Value = 0.004785582017499037
$ws.Cells.Item(3, 20).Value = 0.004785582017499036

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 42.517849
$ws.Cells.Item(4, 8).Value = 127.553547
$ws.Cells.Item(4, 9).Value = 0.02311569285614191
$ws.Cells.Item(4, 10).Value = 0.02311569285614191
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.745494
$ws.Cells.Item(4, 14).Value = 2.236482
$ws.Cells.Item(4, 15).Value = 0.119046603886843
$ws.Cells.Item(4, 16).Value = 0.119046603886843
$ws.Cells.Item(4, 17).Value = 31.696801322406
$ws.Cells.Item(4, 18).Value = 285.271211901654
$ws.Cells.Item(4, 19).Value = 0.002751844731015054
$ws.Cells.Item(4, 20).Value = 0.002751844731015053

$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1644.738728666666
$ws.Cells.Item(5, 8).Value = 4934.216186
$ws.Cells.Item(5, 9).Value = 0.8941956419399297
$ws.Cells.Item(5, 10).Value = 0.8941956419399296
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 4.220261333333333
$ws.Cells.Item(5, 14).Value = 12.660784
$ws.Cells.Item(5, 15).Value = 0.6739259863235564
$ws.Cells.Item(5, 16).Value = 0.6739259863235564
$ws.Cells.Item(5, 17).Value = 6941.227260027757
$ws.Cells.Item(5, 18).Value = 62471.04534024982
$ws.Cells.Item(5, 19).Value = 0.6026216799605928
$ws.Cells.Item(5, 20).Value = 0.6026216799605927

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1644.738728666666
$ws.Cells.Item(6, 8).Value = 4934.216186
$ws.Cells.Item(6, 9).Value = 0.8941956419399297
$ws.Cells.Item(6, 10).Value = 0.8941956419399296
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 1.296447666666667
$ws.Cells.Item(6, 14).Value = 3.889343
$ws.Cells.Item(6, 15).Value = 0.2070274097896007
$ws.Cells.Item(6, 16).Value = 0.2070274097896007
$ws.Cells.Item(6, 17).Value = 2132.3176870562
$ws.Cells.Item(6, 18).Value = 19190.8591835058
$ws.Cells.Item(6, 19).Value = 0.1851230075959729
$ws.Cells.Item(6, 20).Value = 0.1851230075959729

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1644.738728666666
$ws.Cells.Item(7, 8).Value = 4934.216186
$ws.Cells.Item(7, 9).Value = 0.8941956419399297
$ws.Cells.Item(7, 10).Value = 0.8941956419399296
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.745494
$ws.Cells.Item(7, 14).Value = 2.236482
$ws.Cells.Item(7, 15).Value = 0.119046603886843
$ws.Cells.Item(7, 16).Value = 0.119046603886843
$ws.Cells.Item(7, 17).Value = 1226.142853788628
$ws.Cells.Item(7, 18).Value = 11035.28568409765
$ws.Cells.Item(7, 19).Value = 0.1064509543833641
$ws.Cells.Item(7, 20).Value = 0.1064509543833641

$row8 = @("sCs", "Bgn", "Fgfr3", "ECs", 3, 1, 152.093394, 456.280182, 0.08268866520392831, 0.0826886652039283, 2, 0.6666666666666666, 4.220261333333333, 12.660784, 0.6739259863235564, 0.6739259863235564, 641.8738697536321, 5776.864827782688, 0.05572604025533572, 0.05572604025533572)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i+1).Value = $row8[$i] }

$row9 = @("sCs", "Bgn", "Fgfr3", "FAPs", 3, 1, 152.093394, 456.280182, 0.08268866520392831, 0.0826886652039283, 3, 1, 1.296447666666667, 3.889343, 0.2070274097896007, 0.2070274097896007, 197.181125766714, 1774.630131900426, 0.01711882017612876, 0.01711882017612876)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i+1).Value = $row9[$i] }

$row10 = @("sCs", "Bgn", "Fgfr3", "sCs", 3, 1, 152.093394, 456.280182, 0.08268866520392831, 0.0826886652039283, 3, 1, 0.745494, 2.236482, 0.119046603886843, 0.119046603886843, 113.384712666636, 1020.462413999724, 0.009843804772463834, 0.009843804772463832)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i+1).Value = $row10[$i] }
